$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.339.78"
$ws.Range("E2").Value = "  -2.67%  "
$ws.Range("D3").Value = "1.935.14"
$ws.Range("E3").Value = "  -2.65%  "
$ws.Range("E4").Value = "  -0.85%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7094"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.49%  "
$ws.Range("E7").Value = "  -0.67%  "
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.75"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07297"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8052"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08069"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").Value = "1.935.60"
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.486"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").Value = "30.353.46"
$ws.Range("E17").Value = "  -2.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "253.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008213"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.797"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.87%  "
$ws.Range("D21").Value = "2.191.86"
$ws.Range("E21").Value = "  -2.96%  "
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.996"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.752"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.347"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1289"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.351"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.542"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.422"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.168"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05194"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.262"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7478"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.784"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01968"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.813"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "79.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.428"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4532"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.020"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8449"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.762"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.451"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4175"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06033"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.29%  "
